$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Data"
$ws.Name = "Data"

# Update A1 and add A2 with the new data
$ws.Range("A1").Value = "`"- The Gourmet Bistro"
$ws.Range("A2").Value = "Sunshine Café`""
